$d = $word.ActiveDocument

$rsq = [char]0x2019   # RIGHT SINGLE QUOTATION MARK (curly apostrophe)

# The last paragraph in the doc currently ends with a hidden "_GoBack" bookmark
# (zero-length, right after the run text, just before the paragraph mark). We
# delete it now and re-create it at the end of the new content once everything
# has been typed in, mirroring where Word leaves it after the latest edit.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Anchor off the final paragraph ("Chapter onbetfair terminology and api structure")
$tail = $d.Paragraphs.Last.Range

function New-Para([string]$text) {
    $tail.InsertParagraphAfter()
    $tail.Collapse(0)
    $tail.Move(4, 1) | Out-Null
    if ($text.Length -gt 0) {
        $tail.InsertAfter($text)
        $tail.Collapse(0)
    }
}

# Blank paragraph right after the chapter heading line
New-Para ""

New-Para "Need to feed markets closed and data into analysis"
New-Para "This means more can be predicted"

New-Para ""

New-Para ("General idea is that sports not popular In the uk aren" + $rsq + "t viable because betfair is uk/aus based")
New-Para "All popular American ones are not viable except American football, because season length? General popularity? It is the top us sport in usa and expanding I guess, London has been a focus and there was nfl Europe."
New-Para "Testing and eval -> storing json to reserialize"
New-Para "Program design, each sport needs a special module so you can cdeal with game time conversion etc"

# Final paragraph, built from several runs, two of which ("st") are superscript
# ordinal suffixes.
New-Para "Each module needs to support each market, mostly separately because "
$tail.InsertAfter(("game times are different, football halves, nfl quarters, basketball quarters, and there" + $rsq + "s market specific to those, so say if all 1"))
$tail.Collapse(0)

$supStart = $tail.Start
$tail.InsertAfter("st")
$supEnd = $tail.End
$d.Range($supStart, $supEnd).Font.Superscript = $true
$tail.Collapse(0)

$tail.InsertAfter(" quarter markets end then its lieky that the 1")
$tail.Collapse(0)

$supStart2 = $tail.Start
$tail.InsertAfter("st")
$supEnd2 = $tail.End
$d.Range($supStart2, $supEnd2).Font.Superscript = $true
$tail.Collapse(0)

$tail.InsertAfter(" quarter is done.")
$tail.Collapse(0)

# Re-seat "_GoBack" as a zero-length bookmark at the very end of the paragraph
# we just typed (after the last run, before the paragraph mark). A bookmark
# range collapsed exactly at a paragraph's content end gets normalised to span
# the whole paragraph, so we anchor it on a throwaway trailing character and
# then erase just that character, leaving the bookmark pinned in place.
$tail.InsertAfter("X")
$tail.Collapse(0)
$tail.MoveStart(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $tail) | Out-Null
$tail.Text = ""
